$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot old values for columns D, J, K, L, M, P (rows 2-48)
$oldD = @{}
$oldJ = @{}
$oldK = @{}
$oldL = @{}
$oldM = @{}
$oldP = @{}
for ($r = 2; $r -le 48; $r++) {
    $oldD[$r] = $ws.Cells.Item($r, 4).Value2
    $oldJ[$r] = $ws.Cells.Item($r, 10).Value2
    $oldK[$r] = $ws.Cells.Item($r, 11).Value2
    $oldL[$r] = $ws.Cells.Item($r, 12).Value2
    $oldM[$r] = $ws.Cells.Item($r, 13).Value2
    $oldP[$r] = $ws.Cells.Item($r, 16).Value2
}

# Mapping: new row -> source (old) row, derived from the target diff
$rowMap = @{}
$rowMap[2] = 14
$rowMap[3] = 34
$rowMap[4] = 19
$rowMap[5] = 3
$rowMap[6] = 48
$rowMap[7] = 9
$rowMap[8] = 4
$rowMap[9] = 15
$rowMap[10] = 2
$rowMap[11] = 47
$rowMap[12] = 28
$rowMap[13] = 39
$rowMap[14] = 43
$rowMap[15] = 11
$rowMap[16] = 8
$rowMap[17] = 10
$rowMap[18] = 20
$rowMap[19] = 5
$rowMap[20] = 44
$rowMap[21] = 16
$rowMap[22] = 35
$rowMap[23] = 24
$rowMap[24] = 26
$rowMap[25] = 21
$rowMap[26] = 40
$rowMap[27] = 17
$rowMap[28] = 30
$rowMap[29] = 46
$rowMap[30] = 37
$rowMap[31] = 38
$rowMap[32] = 22
$rowMap[33] = 18
$rowMap[34] = 41
$rowMap[35] = 29
$rowMap[36] = 7
$rowMap[37] = 45
$rowMap[38] = 23
$rowMap[39] = 42
$rowMap[40] = 36
$rowMap[41] = 33
$rowMap[42] = 31
$rowMap[43] = 27
$rowMap[44] = 25
$rowMap[45] = 6
$rowMap[46] = 12
$rowMap[47] = 32
$rowMap[48] = 13

foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    $ws.Cells.Item($newRow, 4).Value = $oldD[$srcRow]
    $ws.Cells.Item($newRow, 10).Value = $oldJ[$srcRow]
    $ws.Cells.Item($newRow, 11).Value = $oldK[$srcRow]
    $ws.Cells.Item($newRow, 12).Value = $oldL[$srcRow]
    $ws.Cells.Item($newRow, 13).Value = $oldM[$srcRow]
    $ws.Cells.Item($newRow, 16).Value = $oldP[$srcRow]
}
